$wb = $excel.ActiveWorkbook

# Sheet 1: 展览 (rId1)
$ws1 = $wb.Worksheets.Item(1)
$ws1.Range("F2").Value = 280
$ws1.Range("F5").Value = 2906
$ws1.Range("F7").Value = 122
$ws1.Range("F9").Value = 1530
$ws1.Range("F11").Value = 829
$ws1.Range("F13").Value = 2602
$ws1.Range("F16").Value = 6563
$ws1.Range("F18").Value = 6228
$ws1.Range("F20").Value = 2169
$ws1.Range("F21").Value = 2995
$ws1.Range("F22").Value = 3416
$ws1.Range("F24").Value = 28
$ws1.Range("F25").Value = 1700
$ws1.Range("F26").Value = 61
$ws1.Range("F28").Value = 856
$ws1.Range("F31").Value = 351
$ws1.Range("F32").Value = 1068
$ws1.Range("F33").Value = 2307
$ws1.Range("F35").Value = 146
$ws1.Range("F36").Value = 328
$ws1.Range("F37").Value = 881
$ws1.Range("F38").Value = 182
$ws1.Range("F39").Value = 417
$ws1.Range("F40").Value = 479

# Sheet 2: 演出 (rId2)
$ws2 = $wb.Worksheets.Item(2)
$ws2.Range("F23").Value = 1

# Sheet 4: 全部类型 (rId4)
$ws4 = $wb.Worksheets.Item(4)
$ws4.Range("F5").Value = 280
$ws4.Range("F9").Value = 2906
$ws4.Range("F12").Value = 1530
$ws4.Range("F14").Value = 829
$ws4.Range("F17").Value = 2602
$ws4.Range("F23").Value = 6564
$ws4.Range("F25").Value = 6228
$ws4.Range("F26").Value = 2169
$ws4.Range("F27").Value = 2995
$ws4.Range("F28").Value = 3416
$ws4.Range("F33").Value = 1700
$ws4.Range("F37").Value = 856
$ws4.Range("F40").Value = 351
$ws4.Range("F42").Value = 2307
$ws4.Range("F44").Value = 146
$ws4.Range("F45").Value = 328
$ws4.Range("F46").Value = 881
$ws4.Range("F47").Value = 182
$ws4.Range("F48").Value = 417
$ws4.Range("F49").Value = 479
